{"js": "// Replace the old (2018 / Perseus) campaign-date paragraphs with the new\n// \"Kampagnendaten Cygnus: 10.-19.\" text, as plain unformatted text in a\n// single run, matching the target diff.\nconst OLD_TEXT =\n  \"Kampagnendaten 2018 f\u00fcr das Sternbild Perseus: 30. Oktober - 8. November und 29. November - 8. Dezember\";\nconst NEW_TEXT = \"Kampagnendaten Cygnus: 10.-19.\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  if (paragraph.text === OLD_TEXT) {\n    // Drop every run/formatting from the paragraph first, then insert the\n    // replacement text as a brand-new, unformatted run (mirrors the diff,\n    // which leaves the new run with no <w:rPr> at all).\n    paragraph.clear();\n    paragraph.insertText(NEW_TEXT, \"Start\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the old (2018 / Perseus) campaign-date paragraphs with the new\n# \"Kampagnendaten Cygnus: 10.-19.\" text, as plain unformatted text in a\n# single run, matching the target diff. There are four occurrences in the\n# document (identical visible text, split across several runs each).\n\n$d = $word.ActiveDocument\n\n$oldText = \"Kampagnendaten 2018 f\u00fcr das Sternbild Perseus: 30. Oktober - 8. November und 29. November - 8. Dezember\"\n$newText = \"Kampagnendaten Cygnus: 10.-19.\"\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $r = $p.Range\n    $t = $r.Text\n\n    # $r.Text includes the trailing paragraph mark, so compare against the\n    # text minus that final character.\n    if ($t.Length -ge 1 -and $t.Substring(0, $t.Length - 1) -eq $oldText) {\n        $paraStart = $r.Start\n        $paraEnd = $r.End\n\n        # Delete just the old run text, leaving the paragraph mark (at\n        # paraEnd - 1) untouched so the paragraph itself is preserved\n        # (deleting all the way through the mark merges with the next\n        # paragraph).\n        $textRange = $d.Range($paraStart, $paraEnd - 1)\n        $textRange.Delete()\n\n        # Insert the replacement as a brand-new run with no inherited\n        # character formatting (mirrors the diff, whose new run carries no\n        # <w:rPr> at all).\n        $insertionPoint = $d.Range($paraStart, $paraStart)\n        $insertionPoint.InsertBefore($newText)\n    }\n}\n"}
